$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:A5").Value = 45108

$ws.Range("C2").Value = 1485
$ws.Range("C3").Value = 154
$ws.Range("C4").Value = 280
$ws.Range("C5").Value = 154

$ws.Range("B14").Select()
